$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New submission row 11 mirrors row 10's data exactly, except for a later
# timestamp in column H. Copy row 10 first so the text-typed / empty-string
# cells (A, B, C, D, E, F, G) keep the same value typing as the rest of the
# sheet, then overwrite just the timestamp cell.
$ws.Range("A10:H10").Copy($ws.Range("A11:H11"))
$ws.Cells.Item(11, 8).Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٥:٣٧:٤٧ م"
